$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build list of (cell, newValue) updates derived from the commit diff.
# NumberFormat is forced to Text ("@") before assignment so that the
# numeric-looking strings (prices, percentages, hour codes) are stored
# as text, matching the original inlineStr cell type instead of being
# auto-converted by Excel into numbers/percentages.
$updates = @(
    @{Cell='D2'; Value='321.71'}
    @{Cell='E2'; Value='-2.69%'}
    @{Cell='G2'; Value='10'}
    @{Cell='D3'; Value='42.64'}
    @{Cell='E3'; Value='-5.97%'}
    @{Cell='G3'; Value='10'}
    @{Cell='D4'; Value='5.207'}
    @{Cell='E4'; Value='-7.11%'}
    @{Cell='G4'; Value='10'}
    @{Cell='D5'; Value='0.08204'}
    @{Cell='E5'; Value='-1.76%'}
    @{Cell='G5'; Value='10'}
    @{Cell='D6'; Value='4.301'}
    @{Cell='E6'; Value='-3.32%'}
    @{Cell='G6'; Value='10'}
    @{Cell='D7'; Value='1.803'}
    @{Cell='E7'; Value='-13.41%'}
    @{Cell='G7'; Value='10'}
    @{Cell='D8'; Value='0.9314'}
    @{Cell='E8'; Value='-4.02%'}
    @{Cell='G8'; Value='10'}
    @{Cell='D9'; Value='0.1111'}
    @{Cell='E9'; Value='-5.85%'}
    @{Cell='G9'; Value='10'}
    @{Cell='D10'; Value='0.1865'}
    @{Cell='E10'; Value='-2.96%'}
    @{Cell='G10'; Value='10'}
    @{Cell='D11'; Value='0.09456'}
    @{Cell='E11'; Value='-3.95%'}
    @{Cell='G11'; Value='10'}
    @{Cell='D12'; Value='0.04683'}
    @{Cell='E12'; Value='0.12%'}
    @{Cell='G12'; Value='10'}
    @{Cell='D13'; Value='7.424'}
    @{Cell='E13'; Value='-28.27%'}
    @{Cell='G13'; Value='10'}
    @{Cell='D14'; Value='0.1059'}
    @{Cell='E14'; Value='-0.22%'}
    @{Cell='G14'; Value='10'}
    @{Cell='D15'; Value='0.001301'}
    @{Cell='E15'; Value='1.61%'}
    @{Cell='G15'; Value='10'}
    @{Cell='D16'; Value='0.005714'}
    @{Cell='E16'; Value='-5.81%'}
    @{Cell='G16'; Value='10'}
    @{Cell='D17'; Value='3.363'}
    @{Cell='E17'; Value='-0.32%'}
    @{Cell='G17'; Value='10'}
    @{Cell='E18'; Value='-0.87%'}
    @{Cell='G18'; Value='10'}
    @{Cell='D19'; Value='0.3377'}
    @{Cell='E19'; Value='1.03%'}
    @{Cell='G19'; Value='10'}
    @{Cell='D20'; Value='0.1387'}
    @{Cell='E20'; Value='-0.33%'}
    @{Cell='G20'; Value='10'}
    @{Cell='D21'; Value='0.2546'}
    @{Cell='E21'; Value='-3.43%'}
    @{Cell='G21'; Value='10'}
    @{Cell='D22'; Value='0.04150'}
    @{Cell='E22'; Value='-0.52%'}
    @{Cell='G22'; Value='10'}
    @{Cell='D23'; Value='0.001245'}
    @{Cell='E23'; Value='-5.09%'}
    @{Cell='G23'; Value='10'}
    @{Cell='D24'; Value='0.004353'}
    @{Cell='E24'; Value='-4.50%'}
    @{Cell='G24'; Value='10'}
    @{Cell='E25'; Value='-7.92%'}
    @{Cell='G25'; Value='10'}
    @{Cell='D26'; Value='0.0002981'}
    @{Cell='E26'; Value='-20.47%'}
    @{Cell='G26'; Value='10'}
    @{Cell='G27'; Value='10'}
    @{Cell='G28'; Value='10'}
    @{Cell='G29'; Value='10'}
    @{Cell='G30'; Value='10'}
    @{Cell='G31'; Value='10'}
    @{Cell='G32'; Value='10'}
    @{Cell='G33'; Value='10'}
    @{Cell='G34'; Value='10'}
    @{Cell='G35'; Value='10'}
    @{Cell='G36'; Value='10'}
    @{Cell='G37'; Value='10'}
    @{Cell='D38'; Value='0.02766'}
    @{Cell='E38'; Value='2.21%'}
    @{Cell='G38'; Value='10'}
    @{Cell='D39'; Value='0.05611'}
    @{Cell='E39'; Value='-2.60%'}
    @{Cell='G39'; Value='10'}
    @{Cell='D40'; Value='0.008046'}
    @{Cell='E40'; Value='2.50%'}
    @{Cell='G40'; Value='10'}
    @{Cell='E41'; Value='-2.39%'}
    @{Cell='G41'; Value='10'}
    @{Cell='D42'; Value='0.006549'}
    @{Cell='E42'; Value='-10.29%'}
    @{Cell='G42'; Value='10'}
    @{Cell='D43'; Value='0.002093'}
    @{Cell='E43'; Value='3.26%'}
    @{Cell='G43'; Value='10'}
    @{Cell='D44'; Value='0.008342'}
    @{Cell='E44'; Value='-8.81%'}
    @{Cell='G44'; Value='10'}
    @{Cell='D45'; Value='0.3495'}
    @{Cell='E45'; Value='-1.45%'}
    @{Cell='G45'; Value='10'}
    @{Cell='D46'; Value='0.00006948'}
    @{Cell='E46'; Value='-2.28%'}
    @{Cell='G46'; Value='10'}
    @{Cell='D47'; Value='0.00000000750'}
    @{Cell='E47'; Value='-0.15%'}
    @{Cell='G47'; Value='10'}
    @{Cell='D48'; Value='0.003477'}
    @{Cell='E48'; Value='-0.63%'}
    @{Cell='G48'; Value='10'}
    @{Cell='D49'; Value='0.003532'}
    @{Cell='E49'; Value='0.76%'}
    @{Cell='G49'; Value='10'}
    @{Cell='D50'; Value='0.00002101'}
    @{Cell='E50'; Value='-0.15%'}
    @{Cell='G50'; Value='10'}
    @{Cell='D51'; Value='0.0002001'}
    @{Cell='E51'; Value='-0.15%'}
    @{Cell='G51'; Value='10'}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
}
